$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for rows 2-27
# from serial 45175 (2023-09-06) to serial 45183 (2023-09-14)
for ($row = 2; $row -le 27; $row++) {
    $ws.Cells.Item($row, 3).Value = 45183
}
